# Swap the full data (columns B:AD) between pairs of rows: (130,131), (132,133),
# (251,252), (282,283), (298,299). Column A (row index) stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Parallel arrays (avoid nested-array foreach quirks): row1s[i] swaps with row2s[i].
$row1s = @(130, 132, 251, 282, 298)
$row2s = @(131, 133, 252, 283, 299)

for ($i = 0; $i -lt $row1s.Length; $i++) {
    $r1 = $row1s[$i]
    $r2 = $row2s[$i]

    $rng1 = $ws.Range("B${r1}:AD${r1}")
    $rng2 = $ws.Range("B${r2}:AD${r2}")

    $v1 = $rng1.Value2
    $v2 = $rng2.Value2

    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}
